# Applies the commit: insert two new daily-price records (rows) into the
# "Limón" sheet, right before the existing row that used to be row 632
# (Fecha 44467 / "1a amarillo"). Everything from the old row 632 onward
# shifts down by two rows, and the two new rows are populated with their
# own data (Fecha 44946, calidades "1a plateado" / "2a plateado").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 632-633; this shifts rows 632:699 down to 634:701
# and extends the sheet dimension to A1:T701 automatically.
$ws.Rows("632:633").Insert()

# --- New row 632 ---
$ws.Cells.Item(632, 1).Value = 4
$ws.Cells.Item(632, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(632, 3).Value = "Los Lagos"
$ws.Cells.Item(632, 4).Value = 44946
$ws.Cells.Item(632, 5).Value = 10
$ws.Cells.Item(632, 6).Value = "Fruta"
$ws.Cells.Item(632, 7).Value = 100102
$ws.Cells.Item(632, 8).Value = "Cítricos"
$ws.Cells.Item(632, 9).Value = 100102003
$ws.Cells.Item(632, 10).Value = "Limón"
$ws.Cells.Item(632, 11).Value = "Sin especificar"
$ws.Cells.Item(632, 12).Value = "1a plateado"
$ws.Cells.Item(632, 13).Value = 1200
$ws.Cells.Item(632, 14).Value = 20000
$ws.Cells.Item(632, 15).Value = 21000
$ws.Cells.Item(632, 16).Value = 20500
$ws.Cells.Item(632, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(632, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(632, 19).Value = 1139
$ws.Cells.Item(632, 20).Value = 18

# --- New row 633 ---
$ws.Cells.Item(633, 1).Value = 4
$ws.Cells.Item(633, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(633, 3).Value = "Los Lagos"
$ws.Cells.Item(633, 4).Value = 44946
$ws.Cells.Item(633, 5).Value = 10
$ws.Cells.Item(633, 6).Value = "Fruta"
$ws.Cells.Item(633, 7).Value = 100102
$ws.Cells.Item(633, 8).Value = "Cítricos"
$ws.Cells.Item(633, 9).Value = 100102003
$ws.Cells.Item(633, 10).Value = "Limón"
$ws.Cells.Item(633, 11).Value = "Sin especificar"
$ws.Cells.Item(633, 12).Value = "2a plateado"
$ws.Cells.Item(633, 13).Value = 600
$ws.Cells.Item(633, 14).Value = 18000
$ws.Cells.Item(633, 15).Value = 18000
$ws.Cells.Item(633, 16).Value = 18000
$ws.Cells.Item(633, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(633, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(633, 19).Value = 1000
$ws.Cells.Item(633, 20).Value = 18
